# Auto-generated Excel COM-interop script
# Applies scheduled-runner market data refresh values to Sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2728.8572
$ws.Range("J32").Value = 2728.8572
$ws.Range("L32").Value = 2728.8572
$ws.Range("N32").Value = -3380.8572
$ws.Range("H33").Value = 234.66667
$ws.Range("I33").Value = 129.77777
$ws.Range("K33").Value = 129.77777
$ws.Range("M33").Value = 99.22223
$ws.Range("H40").Value = 2800
$ws.Range("I40").Value = 3033.3333
$ws.Range("K40").Value = 3033.3333
$ws.Range("M40").Value = -2858.3333
$ws.Range("H100").Value = 1945.0555
$ws.Range("I100").Value = 1751.25
$ws.Range("K100").Value = 1751.25
$ws.Range("M100").Value = -1210.25
$ws.Range("H137").Value = 1107.983
$ws.Range("I137").Value = 848.1724
$ws.Range("J137").Value = 1359.1333
$ws.Range("K137").Value = 2544.5172
$ws.Range("L137").Value = 4077.3999
$ws.Range("M137").Value = 5.48279999999977
$ws.Range("N137").Value = -9177.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3527.1587
$ws.Range("I32").Value = 3161.85
$ws.Range("K32").Value = 3161.85
$ws.Range("M32").Value = -2874.85
$ws.Range("H74").Value = 1305.0769
$ws.Range("I74").Value = 814.1
$ws.Range("J74").Value = 2941.6667
$ws.Range("K74").Value = 814.1
$ws.Range("L74").Value = 2941.6667
$ws.Range("M74").Value = 59.89999999999998
$ws.Range("N74").Value = -4689.6667
$ws.Range("H77").Value = 1305.0769
$ws.Range("I77").Value = 814.1
$ws.Range("J77").Value = 2941.6667
$ws.Range("K77").Value = 4070.5
$ws.Range("L77").Value = 14708.3335
$ws.Range("M77").Value = 297.5
$ws.Range("N77").Value = -23444.3335
$ws.Range("H86").Value = 25000
$ws.Range("J86").Value = 25000
$ws.Range("L86").Value = 25000
$ws.Range("N86").Value = -27372
$ws.Range("H89").Value = 25000
$ws.Range("J89").Value = 25000
$ws.Range("L89").Value = 75000
$ws.Range("N89").Value = -86856
$ws.Range("H97").Value = 633.4666999999999
$ws.Range("I97").Value = 508.5
$ws.Range("J97").Value = 1133.3334
$ws.Range("K97").Value = 508.5
$ws.Range("L97").Value = 1133.3334
$ws.Range("M97").Value = -12.5
$ws.Range("N97").Value = -2125.3334
$ws.Range("H102").Value = 27779396
$ws.Range("I102").Value = 27779396
$ws.Range("K102").Value = 27779396
$ws.Range("M102").Value = -27777774
$ws.Range("H110").Value = 1218.1818
$ws.Range("I110").Value = 1003.46155
$ws.Range("J110").Value = 1528.3334
$ws.Range("K110").Value = 1003.46155
$ws.Range("L110").Value = 1528.3334
$ws.Range("M110").Value = 1041.53845
$ws.Range("N110").Value = -5618.3334
$ws.Range("H132").Value = 1232.9656
$ws.Range("I132").Value = 984
$ws.Range("K132").Value = 2952
$ws.Range("M132").Value = -422

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 190.25
$ws.Range("I11").Value = 190.25
$ws.Range("K11").Value = 190.25
$ws.Range("M11").Value = -50.25
$ws.Range("H94").Value = 22728212
$ws.Range("I94").Value = 22728212
$ws.Range("K94").Value = 22728212
$ws.Range("M94").Value = -22727761
$ws.Range("H134").Value = 5906.393
$ws.Range("I134").Value = 1356.6666
$ws.Range("K134").Value = 4069.9998
$ws.Range("M134").Value = -1534.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1724.2222
$ws.Range("I31").Value = 1572
$ws.Range("J31").Value = 2257
$ws.Range("K31").Value = 1572
$ws.Range("L31").Value = 2257
$ws.Range("M31").Value = -1277
$ws.Range("N31").Value = -2847
$ws.Range("H34").Value = 1724.2222
$ws.Range("I34").Value = 1572
$ws.Range("J34").Value = 2257
$ws.Range("K34").Value = 1572
$ws.Range("L34").Value = 2257
$ws.Range("M34").Value = -1370
$ws.Range("N34").Value = -2661
$ws.Range("H62").Value = 16669117
$ws.Range("I62").Value = 2440
$ws.Range("J62").Value = 100002500
$ws.Range("K62").Value = 2440
$ws.Range("L62").Value = 100002500
$ws.Range("M62").Value = -1816
$ws.Range("N62").Value = -100003748
$ws.Range("H65").Value = 16669117
$ws.Range("I65").Value = 2440
$ws.Range("J65").Value = 100002500
$ws.Range("K65").Value = 12200
$ws.Range("L65").Value = 500012500
$ws.Range("M65").Value = -9080
$ws.Range("N65").Value = -500018740
$ws.Range("H134").Value = 890.32556
$ws.Range("I134").Value = 906.5294
$ws.Range("J134").Value = 829.1111
$ws.Range("K134").Value = 2719.5882
$ws.Range("L134").Value = 2487.3333
$ws.Range("M134").Value = -184.5882000000001
$ws.Range("N134").Value = -7557.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H38").Value = 47.090908
$ws.Range("I38").Value = 43.625
$ws.Range("J38").Value = 56.333332
$ws.Range("K38").Value = 130.875
$ws.Range("L38").Value = 168.999996
$ws.Range("M38").Value = 216.125
$ws.Range("N38").Value = -862.999996
$ws.Range("H131").Value = 20409352
$ws.Range("J131").Value = 1580.2858
$ws.Range("L131").Value = 4740.857400000001
$ws.Range("N131").Value = -14820.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1940.1428
$ws.Range("I132").Value = 1405.1
$ws.Range("J132").Value = 3277.75
$ws.Range("K132").Value = 4215.299999999999
$ws.Range("L132").Value = 9833.25
$ws.Range("M132").Value = -1685.299999999999
$ws.Range("N132").Value = -14893.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1234
$ws.Range("J22").Value = 1451
$ws.Range("L22").Value = 1451
$ws.Range("N22").Value = -2041
$ws.Range("H27").Value = 1234
$ws.Range("J27").Value = 1451
$ws.Range("L27").Value = 1451
$ws.Range("N27").Value = -1665
$ws.Range("H46").Value = 1716.6666
$ws.Range("J46").Value = 2075
$ws.Range("L46").Value = 2075
$ws.Range("N46").Value = -2451
$ws.Range("H61").Value = 2315.3333
$ws.Range("I61").Value = 1996.75
$ws.Range("J61").Value = 2952.5
$ws.Range("K61").Value = 1996.75
$ws.Range("L61").Value = 2952.5
$ws.Range("M61").Value = -1794.75
$ws.Range("N61").Value = -3356.5
$ws.Range("H100").Value = 631.9697
$ws.Range("I100").Value = 438.38095
$ws.Range("K100").Value = 438.38095
$ws.Range("M100").Value = 102.61905
$ws.Range("H113").Value = 2315.3333
$ws.Range("I113").Value = 1996.75
$ws.Range("J113").Value = 2952.5
$ws.Range("K113").Value = 1996.75
$ws.Range("L113").Value = 2952.5
$ws.Range("M113").Value = 173.25
$ws.Range("N113").Value = -7292.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 400.33334
$ws.Range("I81").Value = 450.5
$ws.Range("K81").Value = 901
$ws.Range("M81").Value = 160
$ws.Range("H84").Value = 400.33334
$ws.Range("I84").Value = 450.5
$ws.Range("K84").Value = 4505
$ws.Range("M84").Value = 799
$ws.Range("H107").Value = 427.55
$ws.Range("I107").Value = 397.11765
$ws.Range("K107").Value = 1191.35295
$ws.Range("M107").Value = 728.64705

Write-Output "Applied market data refresh: 181 cells set, 1 cells cleared."